$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Widen column F to fit the new, longer notes text.
$ws.Columns.Item(6).ColumnWidth = 44.5703125

# --- Row 17 (Jersey - Implementation) --------------------------------
# Add a "Completed On" date, matching the date format already used in
# column E (copy format from E13 which already carries style s="2").
$ws.Range("E13").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("E17").Value2 = 42648
$ws.Range("F17").Value2 = "jersey.config.server.provider.packages"

# --- Row 19 (JAX-RS with Spring) --------------------------------------
$ws.Rows.Item(19).RowHeight = 30
$ws.Range("E13").Copy()
$ws.Range("E19").PasteSpecial(-4122)
$ws.Range("E19").Value2 = 42650
$ws.Range("F19").Value2 = "Started and completed  10/07/2016 - One hour"

# --- New rows 30 & 31 --------------------------------------------------
$ws.Rows.Item(30).RowHeight = 30
$ws.Range("A30").Value2 = "Creating Login page for note pad application and add roles in db and add descriptor"
$ws.Range("A31").Value2 = " "

$ws.Range("B27").Select()

$excel.CutCopyMode = $false
